$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the header string for B1 ("value" -> "first_release_value")
$ws.Range("B1").Value = "first_release_value"

# Year-end dates (A column) and the corresponding first-release YoY values (B column)
$dates  = @(38717,39082,39447,39813,40178,40543,40908,41274,41639,42004,42369,42735,43100,43465,43830,44196,44561,44926,45291,45657,46022)
$values = @($null,11.51866285751828,6.007290114363029,4.974642158654663,-8.854978371166311,13.37451445936573,8.397157288535361,3.335027872728791,0.1226029201931267,6.350193621343236,5.122443676600863,3.289893304242164,5.201177892156705,3.567305512643082,2.482498593966143,-8.834100858716409,2.402478842946154,0.9490257960172555,-2.275419501954867,-0.139401726460564,$null)

# Extend the date column (A3:A22) with the same bold/centered/bordered
# date style already used on A2, by cloning its formatting.
$ws.Range("A2").Copy()
$ws.Range("A3:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]

    if ($null -ne $values[$i]) {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    } else {
        $ws.Cells.Item($row, 2).Value = $null
    }
}
